$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registro_Errores")

# Widen column D to accommodate longer descriptions.
# Excel's COM ColumnWidth (characters, Calibri 11 default font) maps to a
# slightly different stored "width" units value in the XML; 49.17 rounds to
# the target stored width of 50.
$ws.Columns.Item(4).ColumnWidth = 49.17

$rows = @(
    @(2, "2025-11-13 00:22:47", "sintoma_resuelto", "Virus ransomware eliminado de x_virus.exe", "No", "No", "No", 0),
    @(2, "2025-11-13 00:22:47", "archivo_infectado_detectado", "Archivo spy_tool.exe puesto en cuarentena (Virus: spyware)", "No", "No", "No", 0),
    @(2, "2025-11-13 00:22:47", "archivo_limpio_eliminado", "Error: kernel32.dll era un archivo limpio y fue eliminado", "No", "No", "No", 0),
    @(2, "2025-11-13 00:22:47", "archivo_infectado_detectado", "Archivo adware_bundle.exe puesto en cuarentena (Virus: adware)", "No", "No", "No", 0),
    @(2, "2025-11-13 00:22:47", "archivo_limpio_cuarentena", "Falso positivo: logfile.log era seguro pero fue puesto en cuarentena", "No", "No", "No", 0)
)

$startRow = 8
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
